# Auto-generated Excel COM-interop script
# Adds 2024-06-14 crime-count data: updates column K (2024 totals) and a handful
# of corrected prior-year cells (columns G/J) across the Citywide, By-Neighborhood,
# and individual neighborhood sheets, per the commit 'Add data for 2024-06-14'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3441
$ws.Range("K3").Value = 3415
$ws.Range("G4").Value = 1485
$ws.Range("J4").Value = 1819
$ws.Range("K4").Value = 714
$ws.Range("K6").Value = 4028
$ws.Range("G7").Value = 24711
$ws.Range("J7").Value = 29290
$ws.Range("K7").Value = 11820

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 93
$ws.Range("K8").Value = 790
$ws.Range("K11").Value = 247
$ws.Range("K15").Value = 119
$ws.Range("K19").Value = 370
$ws.Range("K20").Value = 271
$ws.Range("K27").Value = 121
$ws.Range("K29").Value = 625
$ws.Range("K31").Value = 129
$ws.Range("K33").Value = 477
$ws.Range("K34").Value = 56
$ws.Range("K35").Value = 17
$ws.Range("K36").Value = 141
$ws.Range("K41").Value = 103
$ws.Range("G43").Value = 152
$ws.Range("J43").Value = 251
$ws.Range("K43").Value = 107
$ws.Range("K44").Value = 110
$ws.Range("K48").Value = 146
$ws.Range("K49").Value = 70
$ws.Range("K50").Value = 68
$ws.Range("K54").Value = 233
$ws.Range("K55").Value = 126
$ws.Range("K60").Value = 73
$ws.Range("K65").Value = 287
$ws.Range("K67").Value = 457
$ws.Range("K71").Value = 36
$ws.Range("K72").Value = 56
$ws.Range("K75").Value = 41
$ws.Range("K76").Value = 180
$ws.Range("K78").Value = 151
$ws.Range("K79").Value = 305
$ws.Range("K80").Value = 42
$ws.Range("K81").Value = 10
$ws.Range("K83").Value = 254
$ws.Range("K85").Value = 545
$ws.Range("K88").Value = 140
$ws.Range("K89").Value = 159
$ws.Range("K90").Value = 104
$ws.Range("K91").Value = 122
$ws.Range("J94").Value = 330
$ws.Range("K95").Value = 194
$ws.Range("K97").Value = 102
$ws.Range("K99").Value = 204
$ws.Range("G101").Value = 24711
$ws.Range("J101").Value = 29290
$ws.Range("K101").Value = 11820

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 75
$ws.Range("K4").Value = 13
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 186
$ws.Range("K7").Value = 545

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 230
$ws.Range("K6").Value = 261
$ws.Range("K7").Value = 790

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 92
$ws.Range("K3").Value = 82
$ws.Range("K7").Value = 254

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 131
$ws.Range("K3").Value = 176
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 477

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 204

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 43
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 457

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 70
$ws.Range("K6").Value = 108
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 178
$ws.Range("K3").Value = 209
$ws.Range("K7").Value = 625

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 99
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 370

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 94
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J4").Value = 26
$ws.Range("J7").Value = 330

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 38
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 16
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K2").Value = 37
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("G4").Value = 13
$ws.Range("J4").Value = 20
$ws.Range("K6").Value = 46
$ws.Range("G7").Value = 152
$ws.Range("J7").Value = 251
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 10
